$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns that differ between the swapped rows: A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# Swap row 28 with row 29
Swap-Rows 28 29

# Swap row 30 with row 31
Swap-Rows 30 31
